$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: pixel_icons_by_oceansdream.png pickup asset credit
$ws.Range("A7").Value = "pixel_icons_by_oceansdream.png"
$ws.Range("B7").Value = "https://opengameart.org/content/various-inventory-24-pixel-icon-set"
$ws.Range("C7").Value = "CC-BY 3.0, CC-BY-SA 3.0"

# Apply hyperlinks (in the same order the diff introduces new relationship ids)
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.fontspace.com/a-area-kilometer-50-font-f53888") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://freesound.org/people/Whiprealgood/sounds/87535/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://freesound.org/people/suntemple/sounds/253172/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://opengameart.org/content/simple-explosion-bleeds-game-art") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://opengameart.org/content/various-inventory-24-pixel-icon-set") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "https://free-game-assets.itch.io/night-city-street-2d-background-tiles") | Out-Null

# Make sure all newly/previously linked cells use the built-in Hyperlink style
foreach ($addr in @("B2","B3","B4","B5","B7","B9","B10")) {
    $ws.Range($addr).Style = "Hyperlink"
}

# Restore the selected cell as recorded in the saved workbook
$ws.Range("C14").Select() | Out-Null

Write-Host "done"
